# Update countries & provincias Spain
# Applies the 25-Aug-2020 03:35 data refresh to the "Pais" sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Timestamp banner (A1) ---------------------------------------------
$ws.Range("A1").Value = "Datos actualizados a 25 de Agosto de 2020 a las 03:35"

function Set-Row($Row, $Country, $Total, $Nuevos, $Activos, $Recuperados, $Criticos, $MuertesHoy, $Muertes) {
    if ($Country) {
        $ws.Cells.Item($Row, 1).Value = $Country
    }
    $ws.Cells.Item($Row, 2).Value = $Total
    $ws.Cells.Item($Row, 3).Value = $Nuevos
    $ws.Cells.Item($Row, 4).Value = $Activos
    $ws.Cells.Item($Row, 5).Value = $Recuperados
    $ws.Cells.Item($Row, 6).Value = $Criticos
    $ws.Cells.Item($Row, 7).Value = $MuertesHoy
    $ws.Cells.Item($Row, 8).Value = $Muertes
}

# --- Plain value refreshes (row / country unchanged) --------------------

Set-Row 4   $null 5915630 41484 3217947 2516569 0 510 181114   # Estados Unidos
Set-Row 27  $null 125647  751   111694  4870    0 10  9083     # Canada
Set-Row 142 $null 1798    14    634     1118    0 2   46       # Bahamas
Set-Row 164 $null 892     0     832     45      0 0   15       # Santo Tome y Principe
Set-Row 173 $null 430     0     345     84      0 0   1        # Burundi
Set-Row 188 $null 168     1     150     9       0 0   9        # Bermudas

# --- Re-ranked block: rows 176-182 now sorted by updated "Casos totales" -
# New order: San Martin (Parte Holandesa), Papua Nueva Guinea,
# Islas Turcas y Caicos (unchanged), Polinesia Francesa, Mauricio,
# Isla de Man, Eritrea.

Set-Row 176 "San Martin (Parte Holandesa)" 408 12 147 244 0 0 17
Set-Row 177 "Papua Nueva Guinea"           401 0  232 165 0 0 4

# Row 178 (Islas Turcas y Caicos) stays exactly as-is.

Set-Row 179 "Polinesia Francesa" 372 74 148 224 0 0 0
Set-Row 180 "Mauricio"           347 1  335 2   0 0 10
Set-Row 181 "Isla de Man"        336 0  312 0   0 0 24
Set-Row 182 "Eritrea"            306 0  274 32  0 0 0

# Row 183 (Mongolia) stays exactly as-is.

# --- Swap: Montserrat now ranks above Islas Malvinas (rows 214-215) -----

Set-Row 214 "Montserrat"     13 0 12 0 0 0 1
Set-Row 215 "Islas Malvinas" 13 0 13 0 0 0 0
